$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 15 ("Fruta / hortaliza, semanal"),
# pushing all existing rows 15..134 down by one (to 16..135).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new observation.
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 45168
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100112005
$ws.Range("G15").Value = "Puerro"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 70
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("N15").Value = "`$/paquete 20 unidades"
$ws.Range("O15").Value = "Provincia de Chacabuco"
$ws.Range("P15").Value = 400
$ws.Range("Q15").Value = 20
$ws.Range("R15").Value = "Hortaliza"
